$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value2 = 1.03
$ws.Range("F3").Value2 = 2.66
$ws.Range("I3").Value2 = 2.82
$ws.Range("K3").Value2 = 6.8
$ws.Range("N3").Value2 = 1.69
$ws.Range("P3").Value2 = 1.69
$ws.Range("Q3").Value2 = 1.84
$ws.Range("S3").Value2 = 1.84
$ws.Range("V3").Value2 = 1.54
$ws.Range("AD4").Value2 = 23
$ws.Range("AL4").Value2 = 100
$ws.Range("AN4").Value2 = 60
$ws.Range("F4").Value2 = 2.32
$ws.Range("G4").Value2 = 2.58
$ws.Range("H4").Value2 = 3.7
$ws.Range("J4").Value2 = 2.74
$ws.Range("N4").Value2 = 2.22
$ws.Range("O4").Value2 = 1.67
$ws.Range("P4").Value2 = 1.41
$ws.Range("R4").Value2 = 1.14
$ws.Range("T4").Value2 = 2.28
$ws.Range("U4").Value2 = 1.63
$ws.Range("W4").Value2 = 1.64
$ws.Range("Y4").Value2 = 11.5
$ws.Range("Z4").Value2 = 29
$ws.Range("L5").Value2 = 1.47
$ws.Range("O5").Value2 = 1.44
$ws.Range("AF6").Value2 = 17
$ws.Range("I6").Value2 = 3.35
$ws.Range("J6").Value2 = 3.75
$ws.Range("P6").Value2 = 2.16
$ws.Range("R6").Value2 = 1.43
$ws.Range("V6").Value2 = 1.43
$ws.Range("W6").Value2 = 1.71
$ws.Range("AA7").Value2 = 580
$ws.Range("AB7").Value2 = 16
$ws.Range("AC7").Value2 = 16
$ws.Range("AE7").Value2 = 270
$ws.Range("AG7").Value2 = 16.5
$ws.Range("AI7").Value2 = 190
$ws.Range("AJ7").Value2 = 9
$ws.Range("AK7").Value2 = 18.5
$ws.Range("AL7").Value2 = 60
$ws.Range("AM7").Value2 = 190
$ws.Range("AN7").Value2 = 4.2
$ws.Range("AO7").Value2 = 370
$ws.Range("G7").Value2 = 1.32
$ws.Range("I7").Value2 = 13.5
$ws.Range("P7").Value2 = 2.62
$ws.Range("Q7").Value2 = 1.54
$ws.Range("R7").Value2 = 1.7
$ws.Range("T7").Value2 = 2.14
$ws.Range("U7").Value2 = 1.76
$ws.Range("W7").Value2 = 4.1
$ws.Range("X7").Value2 = 26
$ws.Range("AI8").Value2 = 130
$ws.Range("J8").Value2 = 4
$ws.Range("N8").Value2 = 2.94
$ws.Range("Q8").Value2 = 1.94
$ws.Range("W8").Value2 = 2.46
$ws.Range("F9").Value2 = 1.96
$ws.Range("G9").Value2 = 1.97
$ws.Range("K9").Value2 = 3.75
$ws.Range("L9").Value2 = 1.41
$ws.Range("N9").Value2 = 3.65
$ws.Range("R9").Value2 = 1.35
$ws.Range("S9").Value2 = 3.4
$ws.Range("U9").Value2 = 2.04
$ws.Range("V9").Value2 = 1.28
$ws.Range("W9").Value2 = 2.02
$ws.Range("Y9").Value2 = 16.5
$ws.Range("AB10").Value2 = 9
$ws.Range("AF10").Value2 = 8
$ws.Range("AN10").Value2 = 5.2
$ws.Range("I10").Value2 = 14.5
$ws.Range("J10").Value2 = 5.9
$ws.Range("L10").Value2 = 1.3
$ws.Range("N10").Value2 = 5
$ws.Range("T10").Value2 = 2.18
$ws.Range("X10").Value2 = 1000
$ws.Range("Y10").Value2 = 1000
$ws.Range("AA11").Value2 = 120
$ws.Range("AD11").Value2 = 23
$ws.Range("AE11").Value2 = 75
$ws.Range("AG11").Value2 = 11.5
$ws.Range("AI11").Value2 = 1000
$ws.Range("AO11").Value2 = 100
$ws.Range("F11").Value2 = 1.91
$ws.Range("G11").Value2 = 1.97
$ws.Range("I11").Value2 = 4.7
$ws.Range("N11").Value2 = 3
$ws.Range("V11").Value2 = 1.27
$ws.Range("W11").Value2 = 2.02
$ws.Range("Z11").Value2 = 1000
$ws.Range("AM12").Value2 = 120
$ws.Range("F12").Value2 = 2.6
$ws.Range("J12").Value2 = 3.3
$ws.Range("K12").Value2 = 3.4
$ws.Range("M12").Value2 = 1.09
$ws.Range("Q12").Value2 = 2.02
$ws.Range("AN13").Value2 = 80
$ws.Range("G13").Value2 = 3.75
$ws.Range("H13").Value2 = 2.6
$ws.Range("I13").Value2 = 2.86
$ws.Range("N13").Value2 = 1.35
$ws.Range("Q13").Value2 = 1.5
$ws.Range("R13").Value2 = 1.19
$ws.Range("W13").Value2 = 1.36
$ws.Range("Z13").Value2 = 17.5
$ws.Range("AJ14").Value2 = 10.5
$ws.Range("H14").Value2 = 14.5
$ws.Range("L14").Value2 = 1.27
$ws.Range("Q14").Value2 = 1.54
$ws.Range("S14").Value2 = 2.32
$ws.Range("AC15").Value2 = 8.199999999999999
$ws.Range("AG15").Value2 = 12
$ws.Range("F15").Value2 = 1.9
$ws.Range("G15").Value2 = 2
$ws.Range("H15").Value2 = 4.2
$ws.Range("I15").Value2 = 4.8
$ws.Range("J15").Value2 = 3.45
$ws.Range("K15").Value2 = 3.8
$ws.Range("N15").Value2 = 3.3
$ws.Range("P15").Value2 = 1.74
$ws.Range("Q15").Value2 = 2.02
$ws.Range("T15").Value2 = 1.79
$ws.Range("U15").Value2 = 1.92
$ws.Range("V15").Value2 = 1.26
$ws.Range("W15").Value2 = 2
$ws.Range("X15").Value2 = 16
$ws.Range("AH16").Value2 = 18
$ws.Range("AJ16").Value2 = 30
$ws.Range("AM16").Value2 = 100
$ws.Range("AO16").Value2 = 44
$ws.Range("F16").Value2 = 2.22
$ws.Range("P16").Value2 = 1.96
$ws.Range("R16").Value2 = 1.36
$ws.Range("S16").Value2 = 3.5
$ws.Range("U16").Value2 = 2.12
$ws.Range("AA17").Value2 = 390
$ws.Range("AC17").Value2 = 13
$ws.Range("AE17").Value2 = 190
$ws.Range("AF17").Value2 = 9.4
$ws.Range("AN17").Value2 = 8
$ws.Range("F17").Value2 = 1.43
$ws.Range("G17").Value2 = 1.49
$ws.Range("H17").Value2 = 8
$ws.Range("I17").Value2 = 9.4
$ws.Range("P17").Value2 = 1.99
$ws.Range("Q17").Value2 = 1.89
$ws.Range("G18").Value2 = 1.19
$ws.Range("P18").Value2 = 2.74
$ws.Range("G19").Value2 = 3.55
$ws.Range("H19").Value2 = 2.22
$ws.Range("W19").Value2 = 1.39
$ws.Range("X19").Value2 = 16.5
$ws.Range("AA20").Value2 = 490
$ws.Range("AD20").Value2 = 44
$ws.Range("AE20").Value2 = 190
$ws.Range("AG20").Value2 = 11
$ws.Range("AH20").Value2 = 1000
$ws.Range("AJ20").Value2 = 13.5
$ws.Range("AN20").Value2 = 3.95
$ws.Range("F20").Value2 = 1.28
$ws.Range("I20").Value2 = 14
$ws.Range("N20").Value2 = 7
$ws.Range("P20").Value2 = 2.9
$ws.Range("R20").Value2 = 1.76
$ws.Range("S20").Value2 = 2.16
$ws.Range("T20").Value2 = 1.84
$ws.Range("J21").Value2 = 4.4
$ws.Range("N21").Value2 = 5
$ws.Range("Q21").Value2 = 1.59
$ws.Range("T21").Value2 = 1.73
$ws.Range("AN22").Value2 = 3.25
$ws.Range("F22").Value2 = 1.23
$ws.Range("G22").Value2 = 1.24
$ws.Range("K22").Value2 = 8.6
$ws.Range("L22").Value2 = 1.21
$ws.Range("S22").Value2 = 1.87
$ws.Range("T22").Value2 = 1.9
$ws.Range("F23").Value2 = 1.63
$ws.Range("H23").Value2 = 5.3
$ws.Range("Q23").Value2 = 1.66
$ws.Range("U23").Value2 = 2.12
$ws.Range("W23").Value2 = 2.36
